$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assigns $text to the cell at $cellAddr while forcing it to be stored
# as a literal text string (never auto-coerced to a number/date by Excel).
# We do this by writing a formula elsewhere that evaluates to the exact text,
# copying it, and pasting only the *value* into the destination cell. This
# avoids touching the destination cell's number format/style.
function Set-TextValue {
    param($cellAddr, $text)
    $escaped = $text.Replace('"', '""')
    $scratch = $ws.Range("ZZ1000")
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $scratch.Clear()
}

# Plain value assignments (safe: E-column percentages, multi-dot D prices)
$ws.Range("D2").Value = '63.519.54'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '3.254.50'
$ws.Range("E3").Value = '  +3.79%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("D8").Value = '3.249.77'
$ws.Range("E8").Value = '  +4.01%  '
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("E10").Value = '  -0.11%  '
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("E12").Value = '  +0.81%  '
$ws.Range("E13").Value = '  -1.72%  '
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("D15").Value = '3.788.83'
$ws.Range("E15").Value = '  +3.78%  '
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").Value = '3.251.56'
$ws.Range("E17").Value = '  +3.68%  '
$ws.Range("D18").Value = '63.479.50'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("E20").Value = '  -1.00%  '
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("E22").Value = '  +4.14%  '
$ws.Range("E23").Value = '  +4.06%  '
$ws.Range("E24").Value = '  -4.38%  '
$ws.Range("E25").Value = '  +1.46%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("E28").Value = '  +3.97%  '
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("E30").Value = '  +4.41%  '
$ws.Range("E31").Value = '  +2.16%  '
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("E33").Value = '  -3.99%  '
$ws.Range("E34").Value = '  -1.57%  '
$ws.Range("E35").Value = '  -0.94%  '
$ws.Range("E36").Value = '  -1.01%  '
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("E38").Value = '  -2.69%  '
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("E40").Value = '  -1.31%  '
$ws.Range("D41").Value = '2.999.54'
$ws.Range("E41").Value = '  +4.24%  '
$ws.Range("E42").Value = '  +1.05%  '
$ws.Range("E43").Value = '  -2.59%  '
$ws.Range("E44").Value = '  -7.28%  '
$ws.Range("E45").Value = '  +2.57%  '
$ws.Range("E46").Value = '  +1.37%  '
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("E49").Value = '  -1.60%  '
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("E51").Value = '  +1.85%  '

# Text-forced assignments (numeric-looking D-column prices that must remain text)
Set-TextValue "D5" '595.37'
Set-TextValue "D6" '141.29'
Set-TextValue "D9" '0.519'
Set-TextValue "D10" '0.148'
Set-TextValue "D12" '0.467'
Set-TextValue "D13" '0.0000248'
Set-TextValue "D14" '34.37'
Set-TextValue "D19" '6.77'
Set-TextValue "D20" '476.83'
Set-TextValue "D21" '14.23'
Set-TextValue "D22" '0.734'
Set-TextValue "D23" '7.98'
Set-TextValue "D24" '83.70'
Set-TextValue "D25" '13.26'
Set-TextValue "D27" '2.74'
Set-TextValue "D28" '7.21'
Set-TextValue "D29" '8.09'
Set-TextValue "D30" '2.14'
Set-TextValue "D31" '27.76'
Set-TextValue "D34" '2.54'
Set-TextValue "D37" '52.78'
Set-TextValue "D39" '0.0393'
Set-TextValue "D40" '423.32'
Set-TextValue "D42" '8.38'
Set-TextValue "D43" '2.76'
Set-TextValue "D44" '0.110'
Set-TextValue "D46" '2.18'
Set-TextValue "D48" '25.92'
Set-TextValue "D49" '2.33'
Set-TextValue "D51" '122.67'

